$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2489807154707996
$ws.Range("C2").Value = 0.9951336656547281
$ws.Range("D2").Value = 0.3929635865958463
$ws.Range("F2").Value = "Pipeline(steps=[('model', RandomForestRegressor(max_depth=3))])"
$ws.Range("G2").Value = 0.1227843119001288
